$wb = $excel.ActiveWorkbook

# --- Tasks sheet: update plan date ---
$wsTasks = $wb.Worksheets.Item("Tasks")
$wsTasks.Range("C2").Value = "Plan for 02/02/2022"

# --- Subsubs sheet: update research/wavelet plan refs and finish counts ---
$wsSubsubs = $wb.Worksheets.Item("Subsubs")
$wsSubsubs.Range("B3").Value = "R-4.1,2,4"
$wsSubsubs.Range("C3").Value = 4
$wsSubsubs.Range("B4").Value = "D-Wavelet"
$wsSubsubs.Range("C4").Value = 3

# --- Subs sheet: move the stored cursor position ---
$wsSubs = $wb.Worksheets.Item("Subs")
$wsSubs.Range("E4").Select()

# Restore "Subsubs" as the active/visible tab (it was selected before the edit)
$wsSubsubs.Activate()
